$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 144.125
$ws.Range("I2").Value = 129.14285
$ws.Range("J2").Value = 249
$ws.Range("K2").Value = 129.14285
$ws.Range("L2").Value = 249
$ws.Range("M2").Value = -16.14285000000001
$ws.Range("N2").Value = -475
$ws.Range("H28").Value = 33333810
$ws.Range("I28").Value = 34483250
$ws.Range("J28").Value = 98
$ws.Range("K28").Value = 34483250
$ws.Range("L28").Value = 98
$ws.Range("M28").Value = -34482765
$ws.Range("N28").Value = -1068
$ws.Range("H33").Value = 833.6070999999999
$ws.Range("I33").Value = 882.5217
$ws.Range("J33").Value = 608.6
$ws.Range("K33").Value = 882.5217
$ws.Range("L33").Value = 608.6
$ws.Range("M33").Value = -653.5217
$ws.Range("N33").Value = -1066.6
$ws.Range("H42").Value = 767
$ws.Range("H49").Value = 290.5
$ws.Range("I49").Value = 72.5
$ws.Range("J49").Value = 508.5
$ws.Range("K49").Value = 217.5
$ws.Range("L49").Value = 1525.5
$ws.Range("N49").Value = -1797.5
$ws.Range("M49").Value = -81.5
$ws.Range("H55").Value = 489.66666
$ws.Range("I55").Value = 149.71428
$ws.Range("K55").Value = 149.71428
$ws.Range("M55").Value = 64.28572
$ws.Range("H62").Value = 6124.75
$ws.Range("I62").Value = 5499.6665
$ws.Range("K62").Value = 5499.6665
$ws.Range("M62").Value = -4875.6665
$ws.Range("H64").Value = 4064.182
$ws.Range("J64").Value = 4451
$ws.Range("L64").Value = 4451
$ws.Range("N64").Value = -4947
$ws.Range("H65").Value = 6124.75
$ws.Range("I65").Value = 5499.6665
$ws.Range("K65").Value = 27498.3325
$ws.Range("M65").Value = -24378.3325
$ws.Range("H67").Value = 4064.182
$ws.Range("J67").Value = 4451
$ws.Range("L67").Value = 4451
$ws.Range("N67").Value = -6167
$ws.Range("H70").Value = 1155.4667
$ws.Range("I70").Value = 2225.3333
$ws.Range("K70").Value = 6675.999899999999
$ws.Range("M70").Value = -6405.999899999999
$ws.Range("H73").Value = 1155.4667
$ws.Range("I73").Value = 2225.3333
$ws.Range("K73").Value = 6675.999899999999
$ws.Range("M73").Value = -5739.999899999999
$ws.Range("H113").Value = 5287.25
$ws.Range("I113").Value = 4849.75
$ws.Range("J113").Value = 5433.0835
$ws.Range("K113").Value = 4849.75
$ws.Range("L113").Value = 5433.0835
$ws.Range("M113").Value = -1595.75
$ws.Range("N113").Value = -11941.0835
$ws.Range("H125").Value = 1358.3478
$ws.Range("I125").Value = 803.05554
$ws.Range("J125").Value = 3357.4
$ws.Range("K125").Value = 7227.49986
$ws.Range("L125").Value = 30216.6
$ws.Range("M125").Value = -4767.49986
$ws.Range("N125").Value = -35136.60000000001
$ws.Range("H127").Value = 85558.25
$ws.Range("J127").Value = 333333
$ws.Range("L127").Value = 999999
$ws.Range("N127").Value = -1009919
$ws.Range("H129").Value = 1525.1578
$ws.Range("I129").Value = 926.4
$ws.Range("J129").Value = 1739
$ws.Range("K129").Value = 2779.2
$ws.Range("L129").Value = 5217
$ws.Range("M129").Value = 2220.8
$ws.Range("N129").Value = -15217
$ws.Range("H132").Value = 6766.383
$ws.Range("I132").Value = 3745.5
$ws.Range("J132").Value = 24028.572
$ws.Range("K132").Value = 11236.5
$ws.Range("L132").Value = 72085.716
$ws.Range("M132").Value = -8706.5
$ws.Range("N132").Value = -77145.716
$ws.Range("H137").Value = 4995.05
$ws.Range("I137").Value = 5694.1377
$ws.Range("J137").Value = 3152
$ws.Range("K137").Value = 17082.4131
$ws.Range("L137").Value = 9456
$ws.Range("M137").Value = -14532.4131
$ws.Range("N137").Value = -14556
$ws.Range("H141").Value = 14121.723
$ws.Range("I141").Value = 14540.647
$ws.Range("K141").Value = 43621.94100000001
$ws.Range("M141").Value = -38441.94100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 23811030
$ws.Range("I2").Value = 31251336
$ws.Range("J2").Value = 2055.6
$ws.Range("K2").Value = 31251336
$ws.Range("L2").Value = 2055.6
$ws.Range("M2").Value = -31251223
$ws.Range("N2").Value = -2281.6
$ws.Range("H21").Value = 6003.5
$ws.Range("I21").Value = 4671.3335
$ws.Range("K21").Value = 4671.3335
$ws.Range("M21").Value = -4297.3335
$ws.Range("H28").Value = 20574
$ws.Range("H32").Value = 6657.047
$ws.Range("I32").Value = 5355.6772
$ws.Range("K32").Value = 5355.6772
$ws.Range("M32").Value = -5068.6772
$ws.Range("H45").Value = 7468.5
$ws.Range("I45").Value = 7406.9585
$ws.Range("J45").Value = 8207
$ws.Range("K45").Value = 7406.9585
$ws.Range("L45").Value = 8207
$ws.Range("M45").Value = -7029.9585
$ws.Range("N45").Value = -8961
$ws.Range("H61").Value = 3074.7778
$ws.Range("I61").Value = 2801.318
$ws.Range("J61").Value = 6082.8335
$ws.Range("K61").Value = 2801.318
$ws.Range("L61").Value = 6082.8335
$ws.Range("M61").Value = -2589.318
$ws.Range("N61").Value = -6506.8335
$ws.Range("H74").Value = 2055.75
$ws.Range("I74").Value = 1992.2858
$ws.Range("K74").Value = 1992.2858
$ws.Range("M74").Value = -1118.2858
$ws.Range("H77").Value = 2055.75
$ws.Range("I77").Value = 1992.2858
$ws.Range("K77").Value = 9961.429
$ws.Range("M77").Value = -5593.429
$ws.Range("H97").Value = 52685628
$ws.Range("I97").Value = 71430424
$ws.Range("J97").Value = 200183.4
$ws.Range("K97").Value = 71430424
$ws.Range("L97").Value = 200183.4
$ws.Range("M97").Value = -71429928
$ws.Range("N97").Value = -201175.4
$ws.Range("H99").Value = 20574
$ws.Range("H110").Value = 2754.4546
$ws.Range("I110").Value = 2909.9
$ws.Range("K110").Value = 2909.9
$ws.Range("M110").Value = -864.9000000000001
$ws.Range("H116").Value = 23811030
$ws.Range("I116").Value = 31251336
$ws.Range("J116").Value = 2055.6
$ws.Range("K116").Value = 31251336
$ws.Range("L116").Value = 2055.6
$ws.Range("M116").Value = -31249042
$ws.Range("N116").Value = -6643.6
$ws.Range("H122").Value = 13767.667
$ws.Range("I122").Value = 2492.7273
$ws.Range("K122").Value = 7478.1819
$ws.Range("M122").Value = -5028.1819
$ws.Range("H132").Value = 559.24
$ws.Range("I132").Value = 499.17392
$ws.Range("J132").Value = 1250
$ws.Range("K132").Value = 1497.52176
$ws.Range("L132").Value = 3750
$ws.Range("M132").Value = 1032.47824
$ws.Range("N132").Value = -8810
$ws.Range("H133").Value = 152497.33
$ws.Range("J133").Value = 152498.5
$ws.Range("L133").Value = 152498.5
$ws.Range("N133").Value = -157558.5
$ws.Range("H136").Value = 3074.7778
$ws.Range("I136").Value = 2801.318
$ws.Range("J136").Value = 6082.8335
$ws.Range("K136").Value = 8403.954000000002
$ws.Range("L136").Value = 18248.5005
$ws.Range("M136").Value = -5853.954000000002
$ws.Range("N136").Value = -23348.5005
$ws.Range("H140").Value = 120000
$ws.Range("J140").Value = 120000
$ws.Range("L140").Value = 120000
$ws.Range("N140").Value = -130360

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 23811030
$ws.Range("I3").Value = 31251336
$ws.Range("J3").Value = 2055.6
$ws.Range("K3").Value = 31251336
$ws.Range("L3").Value = 2055.6
$ws.Range("M3").Value = -31251222
$ws.Range("N3").Value = -2283.6
$ws.Range("H62").Value = 59999.668
$ws.Range("J62").Value = 59999.668
$ws.Range("N62").Value = -61371.668
$ws.Range("L62").Value = 59999.668
$ws.Range("H65").Value = 59999.668
$ws.Range("J65").Value = 59999.668
$ws.Range("N65").Value = -186863.004
$ws.Range("L65").Value = 179999.004
$ws.Range("H99").Value = 3130.2144
$ws.Range("I99").Value = 732.4286
$ws.Range("K99").Value = 732.4286
$ws.Range("M99").Value = 765.5714
$ws.Range("H105").Value = 1787.4828
$ws.Range("I105").Value = 1720.7037
$ws.Range("J105").Value = 2689
$ws.Range("K105").Value = 1720.7037
$ws.Range("L105").Value = 2689
$ws.Range("M105").Value = 26.29629999999997
$ws.Range("N105").Value = -6183
$ws.Range("H107").Value = 5343.5
$ws.Range("I107").Value = 4837.143
$ws.Range("K107").Value = 4837.143
$ws.Range("M107").Value = -2917.143
$ws.Range("H134").Value = 3917.8914
$ws.Range("I134").Value = 3918.9546
$ws.Range("K134").Value = 11756.8638
$ws.Range("M134").Value = -9221.863799999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 46.857143
$ws.Range("I7").Value = 62.333332
$ws.Range("J7").Value = 8.166667
$ws.Range("K7").Value = 62.333332
$ws.Range("L7").Value = 8.166667
$ws.Range("M7").Value = 50.666668
$ws.Range("N7").Value = -234.166667
$ws.Range("H31").Value = 1811.1578
$ws.Range("I31").Value = 1811.1578
$ws.Range("K31").Value = 1811.1578
$ws.Range("M31").Value = -1516.1578
$ws.Range("H34").Value = 1811.1578
$ws.Range("I34").Value = 1811.1578
$ws.Range("K34").Value = 1811.1578
$ws.Range("M34").Value = -1609.1578
$ws.Range("H58").Value = 1138.4
$ws.Range("I58").Value = 557
$ws.Range("K58").Value = 557
$ws.Range("M58").Value = -354
$ws.Range("H62").Value = 200008350
$ws.Range("I62").Value = 250008450
$ws.Range("J62").Value = 7995
$ws.Range("K62").Value = 250008450
$ws.Range("L62").Value = 7995
$ws.Range("M62").Value = -250007826
$ws.Range("N62").Value = -9243
$ws.Range("H65").Value = 200008350
$ws.Range("I65").Value = 250008450
$ws.Range("J65").Value = 7995
$ws.Range("K65").Value = 1250042250
$ws.Range("L65").Value = 39975
$ws.Range("M65").Value = -1250039130
$ws.Range("N65").Value = -46215
$ws.Range("H99").Value = 12603.963
$ws.Range("I99").Value = 8169.143
$ws.Range("K99").Value = 8169.143
$ws.Range("M99").Value = -6671.143
$ws.Range("H103").Value = 18003
$ws.Range("I103").Value = 18003
$ws.Range("K103").Value = 18003
$ws.Range("M103").Value = -16831
$ws.Range("H126").Value = 12603.963
$ws.Range("I126").Value = 8169.143
$ws.Range("K126").Value = 24507.429
$ws.Range("M126").Value = -22037.429
$ws.Range("H132").Value = 5209.548
$ws.Range("I132").Value = 3095.2812
$ws.Range("J132").Value = 11975.2
$ws.Range("K132").Value = 9285.8436
$ws.Range("L132").Value = 35925.60000000001
$ws.Range("M132").Value = -6755.8436
$ws.Range("N132").Value = -40985.60000000001
$ws.Range("H134").Value = 2065.1304
$ws.Range("I134").Value = 1999.9111
$ws.Range("K134").Value = 5999.7333
$ws.Range("M134").Value = -3464.7333
$ws.Range("H136").Value = 1138.4
$ws.Range("I136").Value = 557
$ws.Range("K136").Value = 1671
$ws.Range("M136").Value = 879

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 31564.031
$ws.Range("I4").Value = 117.34615
$ws.Range("J4").Value = 167833
$ws.Range("K4").Value = 352.03845
$ws.Range("L4").Value = 503499
$ws.Range("M4").Value = -240.03845
$ws.Range("N4").Value = -503723
$ws.Range("H5").Value = 731.0769
$ws.Range("I5").Value = 731.0769
$ws.Range("K5").Value = 2193.2307
$ws.Range("M5").Value = -2081.2307
$ws.Range("H40").Value = 117
$ws.Range("I40").Value = 63
$ws.Range("J40").Value = 225
$ws.Range("K40").Value = 252
$ws.Range("L40").Value = 900
$ws.Range("M40").Value = -183
$ws.Range("N40").Value = -1038
$ws.Range("H82").Value = 13116.5
$ws.Range("I82").Value = 14915
$ws.Range("J82").Value = 11831.857
$ws.Range("K82").Value = 44745
$ws.Range("L82").Value = 35495.571
$ws.Range("M82").Value = -44339
$ws.Range("N82").Value = -36307.571
$ws.Range("H85").Value = 13116.5
$ws.Range("I85").Value = 14915
$ws.Range("J85").Value = 11831.857
$ws.Range("K85").Value = 44745
$ws.Range("L85").Value = 35495.571
$ws.Range("M85").Value = -43341
$ws.Range("N85").Value = -38303.571
$ws.Range("H86").Value = 681.9286
$ws.Range("I86").Value = 606.125
$ws.Range("J86").Value = 783
$ws.Range("K86").Value = 1818.375
$ws.Range("L86").Value = 2349
$ws.Range("M86").Value = -632.375
$ws.Range("N86").Value = -4721
$ws.Range("H87").Value = 6356
$ws.Range("I87").Value = 6356
$ws.Range("K87").Value = 19068
$ws.Range("M87").Value = -17820
$ws.Range("H88").Value = 8499
$ws.Range("J88").Value = 8499
$ws.Range("L88").Value = 25497
$ws.Range("N88").Value = -26353
$ws.Range("H89").Value = 681.9286
$ws.Range("I89").Value = 606.125
$ws.Range("J89").Value = 783
$ws.Range("K89").Value = 5455.125
$ws.Range("L89").Value = 7047
$ws.Range("M89").Value = 472.875
$ws.Range("N89").Value = -18903
$ws.Range("H90").Value = 6356
$ws.Range("I90").Value = 6356
$ws.Range("K90").Value = 57204
$ws.Range("M90").Value = -50964
$ws.Range("H91").Value = 8499
$ws.Range("J91").Value = 8499
$ws.Range("L91").Value = 25497
$ws.Range("N91").Value = -28461
$ws.Range("H118").Value = 5472.5
$ws.Range("I118").Value = 5472.5
$ws.Range("K118").Value = 16417.5
$ws.Range("M118").Value = -15174.5
$ws.Range("H135").Value = 731.0769
$ws.Range("I135").Value = 731.0769
$ws.Range("K135").Value = 6579.6921
$ws.Range("M135").Value = -4044.6921
$ws.Range("H139").Value = 3025.7144
$ws.Range("I139").Value = 2086
$ws.Range("J139").Value = 5375
$ws.Range("K139").Value = 6258
$ws.Range("L139").Value = 16125
$ws.Range("M139").Value = -1118
$ws.Range("N139").Value = -26405

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 40970
$ws.Range("J18").Value = 40970
$ws.Range("L18").Value = 40970
$ws.Range("N18").Value = -41556
$ws.Range("H29").Value = 9499
$ws.Range("J29").Value = 9499
$ws.Range("N29").Value = -10079
$ws.Range("L29").Value = 9499
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()
$ws.Range("H102").Value = 4434.4443
$ws.Range("I102").Value = 4324.037
$ws.Range("J102").Value = 4765.6665
$ws.Range("K102").Value = 4324.037
$ws.Range("L102").Value = 4765.6665
$ws.Range("M102").Value = -2702.037
$ws.Range("N102").Value = -8009.6665
$ws.Range("H122").Value = 2753.6667
$ws.Range("I122").Value = 2216.8948
$ws.Range("K122").Value = 6650.6844
$ws.Range("M122").Value = -4200.6844
$ws.Range("H132").Value = 2231.4285
$ws.Range("I132").Value = 1778.2222
$ws.Range("J132").Value = 3047.2
$ws.Range("K132").Value = 5334.6666
$ws.Range("L132").Value = 9141.599999999999
$ws.Range("M132").Value = -2804.6666
$ws.Range("N132").Value = -14201.6
$ws.Range("H134").Value = 146997.12
$ws.Range("J134").Value = 146997.12
$ws.Range("L134").Value = 440991.36
$ws.Range("N134").Value = -446061.36

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3991.5
$ws.Range("I7").Value = 3991.5
$ws.Range("K7").Value = 3991.5
$ws.Range("M7").Value = -3879.5
$ws.Range("H40").Value = 5511.8213
$ws.Range("I40").Value = 5255.7827
$ws.Range("K40").Value = 5255.7827
$ws.Range("M40").Value = -5119.7827
$ws.Range("H43").Value = 13494.526
$ws.Range("J43").Value = 14949.583
$ws.Range("L43").Value = 14949.583
$ws.Range("N43").Value = -15335.583
$ws.Range("H46").Value = 3010.7058
$ws.Range("I46").Value = 1048.1666
$ws.Range("K46").Value = 1048.1666
$ws.Range("M46").Value = -860.1666
$ws.Range("H61").Value = 15600
$ws.Range("I61").Value = 15600
$ws.Range("K61").Value = 15600
$ws.Range("M61").Value = -15398
$ws.Range("H100").Value = 62502508
$ws.Range("H113").Value = 15600
$ws.Range("I113").Value = 15600
$ws.Range("K113").Value = 15600
$ws.Range("M113").Value = -13430
$ws.Range("H122").Value = 2911
$ws.Range("I122").Value = 2673.3845
$ws.Range("K122").Value = 8020.1535
$ws.Range("M122").Value = -5570.1535
$ws.Range("H126").Value = 3991.5
$ws.Range("I126").Value = 3991.5
$ws.Range("K126").Value = 11974.5
$ws.Range("M126").Value = -9504.5
$ws.Range("H132").Value = 50572.22
$ws.Range("I132").Value = 62179.883
$ws.Range("J132").Value = 17683.834
$ws.Range("K132").Value = 186539.649
$ws.Range("L132").Value = 53051.50199999999
$ws.Range("M132").Value = -184009.649
$ws.Range("N132").Value = -58111.50199999999
$ws.Range("H135").Value = 79333.336
$ws.Range("J135").Value = 79333.336
$ws.Range("L135").Value = 79333.336
$ws.Range("N135").Value = -89473.336

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 6700.4
$ws.Range("I32").Value = 6700.4
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 6700.4
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()
$ws.Range("M32").Value = -6383.4
$ws.Range("H113").Value = 606.4091
$ws.Range("I113").Value = 676.05884
$ws.Range("K113").Value = 2028.17652
$ws.Range("M113").Value = 141.82348
$ws.Range("H122").Value = 5343.5625
$ws.Range("I122").Value = 2621.878
$ws.Range("J122").Value = 21284.857
$ws.Range("K122").Value = 7865.634
$ws.Range("L122").Value = 63854.571
$ws.Range("M122").Value = -5415.634
$ws.Range("N122").Value = -68754.571
$ws.Range("H126").Value = 4035260.2
$ws.Range("I126").Value = 5211489
$ws.Range("K126").Value = 15634467
$ws.Range("M126").Value = -15631997
$ws.Range("H132").Value = 3101.88
$ws.Range("I132").Value = 3336.4736
$ws.Range("K132").Value = 10009.4208
$ws.Range("M132").Value = -7479.4208
$ws.Range("H136").Value = 1510.1282
$ws.Range("I136").Value = 1571.3715
$ws.Range("J136").Value = 974.25
$ws.Range("K136").Value = 4714.1145
$ws.Range("L136").Value = 2922.75
$ws.Range("M136").Value = -2164.1145
$ws.Range("N136").Value = -8022.75
